$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert 4 new rows before the old row 23 ("Result" header), shifting
# the old rows 23-28 down to 27-32 ---
$ws.Rows("23:26").Insert()

# --- Populate the 4 new rows (23-26) by cloning the existing "InputLine /
# OutputLine" rows (21-22) which carry the same column layout & styles ---
$ws.Range("B21:M22").Copy($ws.Range("B23:M24"))
$ws.Range("B21:M22").Copy($ws.Range("B25:M26"))

# Row 24 / 26 keep the thin-bordered "bottom of block" look (ht 15) like row 22
$ws.Rows(24).RowHeight = $ws.Rows(22).RowHeight
$ws.Rows(26).RowHeight = $ws.Rows(22).RowHeight

# Fix up the handful of cells whose style differs from the cloned source row
$ws.Range("J23").Copy()
$ws.Range("I23").PasteSpecial(-4122)
$ws.Range("J24").Copy()
$ws.Range("I24").PasteSpecial(-4122)
$ws.Range("J25").Copy()
$ws.Range("H25").PasteSpecial(-4122)
$ws.Range("J26").Copy()
$ws.Range("H26").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New row content/labels + ConnectionResult bug fix (rows 11/12) + the new
# InputConnector/OutputConnector cells, written in the same order the
# shared-string table records them (keeps sharedStrings.xml byte-identical
# to the authored edit)
$ws.Range("B23").Value = "(14)"
$ws.Range("B24").Value = "(15)"
$ws.Range("H11").Value = "input connector (14)"
$ws.Range("I12").Value = "output connector (15)"
$ws.Range("B25").Value = "(16)"
$ws.Range("B26").Value = "(17)"
$ws.Range("J11").Value = "input connector (16)"
$ws.Range("K12").Value = "output connector (17)"
$ws.Range("K11").Value = "array of knots     (6)"

$ws.Range("I23").Value = "InputConnector"
$ws.Range("J23").Value = "InputConnector"
$ws.Range("I24").Value = "OutputConnector"
$ws.Range("J24").Value = "OutputConnector"
$ws.Range("H25").Value = "InputConnector"
$ws.Range("I25").Value = "InputConnector"
$ws.Range("J25").Value = "InputConnector"
$ws.Range("H26").Value = "OutputConnector"
$ws.Range("I26").Value = "OutputConnector"
$ws.Range("J26").Value = "OutputConnector"

$ws.Range("J12").Value = "array of knots     (6)"

# K11 & J12 upgrade from the plain "false" style to the IoConnector style
# (same as their row's other connector cells)
$ws.Range("H11").Copy()
$ws.Range("K11").PasteSpecial(-4122)
$ws.Range("I12").Copy()
$ws.Range("J12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- cosmetic: selection moved by the author while testing ---
$ws.Range("N6").Select()
